# PCBBom.xlsx update — add new resistor / diagnostic-LED BOM rows (17-30),
# restore the shared "price total" formula on the newly-filled rows, and
# wire up a couple of additional DigiKey hyperlinks (M19, M5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 17 — Resistor for Voltage Division (249k) (0805)
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Resistor for Voltage Division (249k) (0805)"
$ws.Range("C17").Value = "R1-R4"
$ws.Range("E17").Value = "RMCF0805FT249K"
$ws.Range("G17").Value = 0.15
$ws.Range("I17").Value = 4
$ws.Range("K17").Formula = "=G17*I17"
$ws.Range("M17").Value = "https://www.digikey.ca/en/products/detail/stackpole-electronics-inc/RMCF0805FT249K/1760185"

# ---------------------------------------------------------------------
# Row 18 — Resistor for Voltage Division (1.3k) (0805)
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Resistor for Voltage Division (1.3k) (0805)"
$ws.Range("C18").Value = "R7-R9"
$ws.Range("E18").Value = "RNCP0805FTD1K30"
$ws.Range("G18").Value = 0.15
$ws.Range("I18").Value = 3
$ws.Range("K18").Formula = "=G18*I18"
$ws.Range("M18").Value = "https://www.digikey.ca/en/products/detail/stackpole-electronics-inc/RNCP0805FTD1K30/2240232"
$ws.Rows.Item(18).RowHeight = 23.25

# ---------------------------------------------------------------------
# Row 19 — Resistor for Input Op Amp (100ohm) (0603)
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Resistor for Input Op Amp (100ohm) (0603)"
$ws.Range("C19").Value = "R10,R11, R27, R28"
$ws.Range("E19").Value = "ESR03EZPJ101"
$ws.Range("G19").Value = 0.15
$ws.Range("I19").Value = 4
$ws.Range("K19").Formula = "=G19*I19"

# ---------------------------------------------------------------------
# Row 20 — Resistors for CMOS Pullup (22k) (0805)
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Resistors for CMOS Pullup (22k) (0805)"
$ws.Range("C20").Value = "R12, R15,R16,R17,R29"
$ws.Range("E20").Value = "CRG0805F22K"
$ws.Range("G20").Value = 0.15
$ws.Range("I20").Value = 5
$ws.Range("K20").Formula = "=G20*I20"
$ws.Range("M20").Value = "https://www.digikey.ca/en/products/detail/te-connectivity-passive-product/CRG0805F22K/2380871"

# ---------------------------------------------------------------------
# Row 21 — Resistor for Opto
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Resistor for Opto"
$ws.Range("C21").Value = "R13"
$ws.Range("E21").Value = "CRGCQ0603F1K2"
$ws.Range("G21").Value = 0.15
$ws.Range("I21").Value = 1
$ws.Range("K21").Formula = "=G21*I21"
$ws.Range("M21").Value = "https://www.digikey.ca/en/products/detail/te-connectivity-passive-product/CRGCQ0603F1K2/8576291"

# ---------------------------------------------------------------------
# Row 22 — Resistor for CMOS Pullup
# ---------------------------------------------------------------------
$ws.Range("A22").Value = "Resistor for CMOS Pullup"
$ws.Range("C22").Value = "R14,R19"
$ws.Range("E22").Value = "RNCP0603FTD10K0"
$ws.Range("G22").Value = 0.15
$ws.Range("I22").Value = 2
$ws.Range("K22").Formula = "=G22*I22"
$ws.Range("M22").Value = "https://www.digikey.ca/en/products/detail/stackpole-electronics-inc/RNCP0603FTD10K0/2240139"

# ---------------------------------------------------------------------
# Row 23 — Resistor for LED (499Ohm)
# ---------------------------------------------------------------------
$ws.Range("A23").Value = "Resistor for LED (499Ohm)"
$ws.Range("C23").Value = "R18,R20"
$ws.Range("E23").Value = "ERA-3AEB4990V"
$ws.Range("G23").Value = 0.51
$ws.Range("I23").Value = 2
$ws.Range("K23").Formula = "=G23*I23"
$ws.Range("M23").Value = "https://www.digikey.ca/en/products/detail/panasonic-electronic-components/ERA-3AEB4990V/2026721"

# ---------------------------------------------------------------------
# Row 24 — Diagnostic LED's (Red)
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "Diagnostic LED's (Red)"
$ws.Range("C24").Value = "D2"
$ws.Range("E24").Value = "XZM2CRK54WA-8"
$ws.Range("G24").Value = 0.84
$ws.Range("I24").Value = 1
$ws.Range("K24").Formula = "=G24*I24"
$ws.Range("M24").Value = "https://www.digikey.ca/en/products/detail/sunled/XZM2CRK54WA-8/8571166"

# ---------------------------------------------------------------------
# Row 25 — Diagnostic LED's (Green)
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Diagnostic LED's (Green)"
$ws.Range("C25").Value = "D7"
$ws.Range("E25").Value = "XZVG54W-8"
$ws.Range("G25").Value = 0.69
$ws.Range("I25").Value = 1
$ws.Range("K25").Formula = "=G25*I25"
$ws.Range("M25").Value = "https://www.digikey.ca/en/products/detail/sunled/XZVG54W-8/8259041"

# ---------------------------------------------------------------------
# Row 26 — Resistor for Voltage Division(150k) (0805)
# ---------------------------------------------------------------------
$ws.Range("A26").Value = "Resistor for Voltage Division(150k) (0805)"
$ws.Range("E26").Value = "ERJ-P06J154V"
$ws.Range("G26").Value = 0.19
$ws.Range("I26").Value = 3
$ws.Range("K26").Formula = "=I26*G26"
$ws.Range("M26").Value = "https://www.digikey.ca/en/products/detail/panasonic-electronic-components/ERJ-P06J154V/525220"
$ws.Range("M26").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Row 27 — Resistor for Voltage Division(4020) (0805)
# ---------------------------------------------------------------------
$ws.Range("A27").Value = "Resistor for Voltage Division(4020) (0805)"
$ws.Range("E27").Value = "RT0805FRE074K02L"
$ws.Range("G27").Value = 0.15
$ws.Range("I27").Value = 1
$ws.Range("K27").Formula = "=I27*G27"
$ws.Range("M27").Value = "https://www.digikey.ca/en/products/detail/yageo/RT0805FRE074K02L/1079285"

# ---------------------------------------------------------------------
# Row 28 — Resistor for Voltage Division (4.7k) (0805)
# ---------------------------------------------------------------------
$ws.Range("A28").Value = "Resistor for Voltage Division (4.7k) (0805)"
$ws.Range("E28").Value = "CRGCQ0805J4K7"
$ws.Range("G28").Value = 0.15
$ws.Range("I28").Value = 1
$ws.Range("K28").Formula = "=G28*I28"
$ws.Range("M28").Value = "https://www.digikey.ca/en/products/detail/te-connectivity-passive-product/CRGCQ0805J4K7/8576740"

# ---------------------------------------------------------------------
# Row 29 — Resistor for Voltage Division (6.8k) (0805)
# ---------------------------------------------------------------------
$ws.Range("A29").Value = "Resistor for Voltage Division (6.8k) (0805)"
$ws.Range("E29").Value = "WR08X6801FTL"
$ws.Range("G29").Value = 0.15
$ws.Range("I29").Value = 1
$ws.Range("K29").Formula = "=G29*I29"
$ws.Range("M29").Value = "https://www.digikey.ca/en/products/detail/walsin-technology-corporation/WR08X6801FTL/13238733"

# ---------------------------------------------------------------------
# Row 30 — Resistor for Relays (Non Connected) (0603)
# ---------------------------------------------------------------------
$ws.Range("A30").Value = "Resistor for Relays (Non Connected) (0603)"

# ---------------------------------------------------------------------
# New hyperlinks: the freshly added M19 row, and the pre-existing M5 cell
# which previously just showed the raw URL as plain text.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("M19"), "https://www.digikey.ca/en/products/detail/rohm-semiconductor/ESR03EZPJ101/1983452") | Out-Null
$ws.Range("M19").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("M5"), "https://www.digikey.ca/en/products/detail/diodes-incorporated/BSS138W-7-F/814992?utm_campaign=buynow&utm_medium=aggregator&WT.z_cid=ref_findchips_standard&utm_source=findchips") | Out-Null
$ws.Range("M5").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Restore the cursor/selection to where the author left off editing.
# ---------------------------------------------------------------------
$ws.Range("J31").Select() | Out-Null
